# Apply the "spec.xlsx" flockBlockersExtra update:
#  - Sheet1 row 1/2 reworked into a cost/conversion header block
#  - Diamond package pricing rows added (3000 / 10000 diamond tiers, cost row)
#  - "Amount (*1.1 ...)" formula note bumped to "*1.15 ..."
#  - Flock row's starting bird count bumped from 7 to 10
#  - Selection cursor moved to J20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: diamond/currency conversion notes -------------------------------
$ws.Range("A1").Value = "1 dia=1000 go, 5000 go=1 dia"
$ws.Range("C1").Value = "90c/1000 0.09c/dia"
$ws.Range("D1").Value = "`$2/3000 0.07c/dia"
$ws.Range("E1").Value = "`$5/10000 0.05c/dia"
$ws.Range("F1").Value = "90c/no ads"

# --- Row 2: cost row + diamond package amounts (shifted one column right) ---
$ws.Range("A2").Value = "cost"

$ws.Range("B2").Value = 500
$ws.Range("B2").HorizontalAlignment = -4108

$ws.Range("C2").Value = 1500
$ws.Range("C2").HorizontalAlignment = -4108

$ws.Range("D2").Value = 4500
$ws.Range("D2").HorizontalAlignment = -4108

$ws.Range("E2").Value = 13500
$ws.Range("E2").HorizontalAlignment = -4108

$ws.Range("F2").Value = 40500
$ws.Range("F2").HorizontalAlignment = -4108

$ws.Range("G2").Value = 121500
$ws.Range("G2").HorizontalAlignment = -4108

$ws.Range("H2").Value = 364500
$ws.Range("H2").HorizontalAlignment = -4108

# I2 used to be the last numeric tier (2048000, centered); it is now a plain
# "1000 diamonds" label with default (general) alignment.
$ws.Range("I2").Style = "Normal"
$ws.Range("I2").Value = "1000 diamonds"

$ws.Range("J2").Value = "3000 diamonds"
$ws.Range("J2").HorizontalAlignment = -4108

$ws.Range("K2").Value = "10000 diamonds"
$ws.Range("K2").HorizontalAlignment = -4108

# --- Row 10: bump the amount-scaling note from *1.1 to *1.15 ----------------
$ws.Range("G10").Value = "Amount (*1.15 each time up to 3x original amounts, then health x2 (&goldx2.5) until 10x original)"

# --- Row 11 (Flock): starting bird count 7 -> 10 -----------------------------
$ws.Range("E11").Value = 10

# --- Selection cursor -> J20 -------------------------------------------------
$ws.Range("J20").Select()
